$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number run, and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/6/2025  Through  1/12/2025"

# --- Crime statistics table updates (rows 14-31) ---
# Row 14
$ws.Range("N14").Value = -100
$ws.Range("M29").Copy()
$ws.Range("N14").PasteSpecial(-4122)
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("I15").Value = 1
$ws.Range("C16").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("L15").Value = 0
$ws.Range("M29").Copy()
$ws.Range("L15").PasteSpecial(-4122)
$ws.Range("N15").Value = -75
# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -8.333333333333
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = -25
$ws.Range("L16").Value = -57.142857142857
$ws.Range("M16").Value = -62.5
$ws.Range("N16").Value = -93.023255813953
# Row 17
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 15.789473684210
$ws.Range("I17").Value = 7
$ws.Range("J17").Value = 7
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = -12.5
$ws.Range("M17").Value = 133.333333333333
$ws.Range("N17").Value = -69.565217391304
# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 14
$ws.Range("E18").Value = -64.285714285714
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 5
$ws.Range("C16").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("J18").Value = 17
$ws.Range("K18").Value = -70.588235294117
$ws.Range("L18").Value = -66.666666666666
$ws.Range("M18").Value = -61.538461538461
$ws.Range("N18").Value = -89.795918367346
# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = -17.391304347826
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 13
$ws.Range("K19").Value = 7.692307692307
$ws.Range("L19").Value = -46.153846153846
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
# Row 20
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 5
$ws.Range("J20").Value = 6
$ws.Range("K20").Value = -16.666666666666
$ws.Range("L20").Value = 25
$ws.Range("M20").Value = -16.666666666666
$ws.Range("N20").Value = -73.684210526315
# Row 21
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -25.925925925925
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 113
$ws.Range("H21").Value = -15.044247787610
$ws.Range("I21").Value = 35
$ws.Range("J21").Value = 47
$ws.Range("K21").Value = -25.531914893617
$ws.Range("L21").Value = -43.548387096774
$ws.Range("M21").Value = -22.222222222222
$ws.Range("N21").Value = -77.272727272727
# Row 22
$ws.Range("D22").Value = "'0"
$ws.Range("C29").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("C29").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 4
# Row 23
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 8
$ws.Range("H23").Value = -46.666666666666
$ws.Range("I23").Value = 3
$ws.Range("J23").Value = 6
$ws.Range("L23").Value = -66.666666666666
# Row 24
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -17.391304347826
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 26.506024096385
$ws.Range("I24").Value = 28
$ws.Range("J24").Value = 37
$ws.Range("K24").Value = -24.324324324324
$ws.Range("L24").Value = -22.222222222222
$ws.Range("M24").Value = -20
# Row 25
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 20.833333333333
$ws.Range("I25").Value = 8
$ws.Range("J25").Value = 10
$ws.Range("K25").Value = -20
$ws.Range("L25").Value = 300
# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 12.5
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 55
$ws.Range("H26").Value = -52.727272727272
$ws.Range("I26").Value = 9
$ws.Range("C16").Copy()
$ws.Range("I26").PasteSpecial(-4122)
$ws.Range("J26").Value = 24
$ws.Range("K26").Value = -62.5
$ws.Range("L26").Value = -35.714285714285
$ws.Range("M26").Value = -18.181818181818
# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1
$ws.Range("C16").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("I27").Value = 1
$ws.Range("C16").Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("L27").Value = 0
$ws.Range("M29").Copy()
$ws.Range("L27").PasteSpecial(-4122)
# Row 28
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 200
# Row 31
$ws.Range("D31").Value = "'0"
$ws.Range("C29").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "'***.*"
$ws.Range("C29").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("F31").Value = "'0"
$ws.Range("C29").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = -100

$excel.CutCopyMode = $false
